$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need an explicit Text format so Excel
# keeps them as strings (matching the source inline-string cells) instead of
# silently converting them to numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row.
$ws.Range("D2").Value = '27.174.41'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '1.780.96'
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '336.71'
$ws.Range("E5").Value = '  -2.05%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '0.3913'
$ws.Range("E7").Value = '  +2.09%  '
$ws.Range("D8").Value = '0.3417'
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").Value = '47.80'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").Value = '1.183'
$ws.Range("E10").Value = '  -4.34%  '
$ws.Range("D11").Value = '0.07411'
$ws.Range("E11").Value = '  -4.64%  '
$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '21.47'
$ws.Range("E13").Value = '  -3.88%  '
$ws.Range("D14").Value = '6.411'
$ws.Range("D15").Value = '1.778.53'
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").Value = '7.075'
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("E17").Value = '  -3.32%  '
$ws.Range("D18").Value = '0.06663'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").Value = '83.11'
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '17.45'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '6.459'
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").Value = '27.131.33'
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("D24").Value = '12.27'
$ws.Range("E24").Value = '  -6.95%  '
$ws.Range("D25").Value = '2.373'
$ws.Range("E25").Value = '  -3.67%  '
$ws.Range("D26").Value = '21.04'
$ws.Range("E26").Value = '  -5.18%  '
$ws.Range("E27").Value = '  -7.74%  '
$ws.Range("D28").Value = '1.431'
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").Value = '155.55'
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").Value = '1.980.79'
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("D31").Value = '133.75'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '3.976'
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D33").Value = '5.936'
$ws.Range("E33").Value = '  -6.99%  '
$ws.Range("D34").Value = '0.08696'
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = '12.83'
$ws.Range("E35").Value = '  -7.84%  '
$ws.Range("D36").Value = '1.617'
$ws.Range("E36").Value = '  -4.26%  '
$ws.Range("D37").Value = '5.361'
$ws.Range("E37").Value = '  -4.90%  '
$ws.Range("D38").Value = '0.6729'
$ws.Range("E38").Value = '  -4.48%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02342'
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06316'
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("D41").Value = '0.2178'
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").Value = '1.237'
$ws.Range("E42").Value = '  -4.78%  '
$ws.Range("D43").Value = '8.391'
$ws.Range("E43").Value = '  -6.71%  '
$ws.Range("D44").Value = '14.14'
$ws.Range("E44").Value = '  -4.22%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = '0.6340'
$ws.Range("E46").Value = '  -4.20%  '
$ws.Range("D47").Value = '3.842'
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("D48").Value = '2.127'
$ws.Range("E48").Value = '  -2.96%  '
$ws.Range("D49").Value = '130.63'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").Value = '0.07110'
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("D51").Value = '78.48'
$ws.Range("E51").Value = '  -2.83%  '
